$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D price cells to remain text (matches original inlineStr typing),
# otherwise numeric-looking strings like "1.012" get auto-converted to numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '30.658.39'
$ws.Range("E2").Value = '  +0.68%  '
$ws.Range("D3").Value = '2.126.34'
$ws.Range("E3").Value = '  +0.94%  '
$ws.Range("D4").Value = '1.012'
$ws.Range("E4").Value = '  +0.63%  '
$ws.Range("D5").Value = '352.61'
$ws.Range("E5").Value = '  +5.40%  '
$ws.Range("E6").Value = '  +0.58%  '
$ws.Range("D7").Value = '0.5279'
$ws.Range("E7").Value = '  +0.98%  '
$ws.Range("D8").Value = '0.4563'
$ws.Range("E8").Value = '  -0.14%  '
$ws.Range("D9").Value = '53.96'
$ws.Range("E9").Value = '  +1.31%  '
$ws.Range("D10").Value = '0.09107'
$ws.Range("E10").Value = '  +1.75%  '
$ws.Range("D11").Value = '1.184'
$ws.Range("E11").Value = '  +0.69%  '
$ws.Range("D12").Value = '24.62'
$ws.Range("E12").Value = '  +1.51%  '
$ws.Range("D13").Value = '2.134.91'
$ws.Range("E13").Value = '  +1.61%  '
$ws.Range("D14").Value = '6.857'
$ws.Range("E14").Value = '  +0.06%  '
$ws.Range("D15").Value = '8.115'
$ws.Range("E15").Value = '  +0.95%  '
$ws.Range("D16").Value = '102.56'
$ws.Range("E16").Value = '  +6.13%  '
$ws.Range("D17").Value = '0.00001179'
$ws.Range("E17").Value = '  +2.77%  '
$ws.Range("D18").Value = '1.012'
$ws.Range("E18").Value = '  +0.61%  '
$ws.Range("D19").Value = '0.06712'
$ws.Range("E19").Value = '  +0.81%  '
$ws.Range("D20").Value = '19.47'
$ws.Range("E20").Value = '  +1.23%  '
$ws.Range("D21").Value = '1.011'
$ws.Range("E21").Value = '  +0.60%  '
$ws.Range("D22").Value = '6.361'
$ws.Range("E22").Value = '  +0.30%  '
$ws.Range("D23").Value = '30.727.26'
$ws.Range("E23").Value = '  +0.70%  '
$ws.Range("D24").Value = '12.91'
$ws.Range("E24").Value = '  +3.16%  '
$ws.Range("D25").Value = '2.389'
$ws.Range("E25").Value = '  +1.24%  '
$ws.Range("D26").Value = '2.392.48'
$ws.Range("E26").Value = '  +1.84%  '
$ws.Range("D27").Value = '22.52'
$ws.Range("E27").Value = '  +1.03%  '
$ws.Range("D28").Value = '2.566'
$ws.Range("E28").Value = '  +0.94%  '
$ws.Range("D29").Value = '164.64'
$ws.Range("E29").Value = '  +0.88%  '
$ws.Range("D30").Value = '137.03'
$ws.Range("E30").Value = '  +2.96%  '
$ws.Range("D31").Value = '1.200'
$ws.Range("E31").Value = '  -1.49%  '
$ws.Range("D32").Value = '0.1083'
$ws.Range("E32").Value = '  +0.90%  '
$ws.Range("D33").Value = '1.668'
$ws.Range("E33").Value = '  +0.47%  '
$ws.Range("D34").Value = '6.384'
$ws.Range("E34").Value = '  +0.16%  '
$ws.Range("D35").Value = '4.023'
$ws.Range("E35").Value = '  +2.03%  '
$ws.Range("D36").Value = '6.178'
$ws.Range("E36").Value = '  +7.99%  '
$ws.Range("D37").Value = '10.34'
$ws.Range("E37").Value = '  -0.27%  '
$ws.Range("D38").Value = '0.02661'
$ws.Range("E38").Value = '  +3.09%  '
$ws.Range("D39").Value = '0.06890'
$ws.Range("E39").Value = '  +0.71%  '
$ws.Range("D40").Value = '0.2326'
$ws.Range("E40").Value = '  +0.93%  '
$ws.Range("D41").Value = '12.56'
$ws.Range("E41").Value = '  -1.15%  '
$ws.Range("D42").Value = '0.6937'
$ws.Range("E42").Value = '  +0.77%  '
$ws.Range("D43").Value = '1.278'
$ws.Range("E43").Value = '  +1.89%  '
$ws.Range("D44").Value = '14.79'
$ws.Range("E44").Value = '  +5.07%  '
$ws.Range("D45").Value = '2.342'
$ws.Range("E45").Value = '  +0.65%  '
$ws.Range("D46").Value = '0.6465'
$ws.Range("E46").Value = '  +1.20%  '
$ws.Range("D47").Value = '3.773'
$ws.Range("E47").Value = '  +2.90%  '
$ws.Range("E48").Value = '  +7.61%  '
$ws.Range("D49").Value = '1.258'
$ws.Range("E49").Value = '  +0.41%  '
$ws.Range("D50").Value = '83.07'
$ws.Range("E50").Value = '  -0.30%  '
$ws.Range("D51").Value = '0.07302'
$ws.Range("E51").Value = '  +2.36%  '
